$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102 (shifts existing row 102..220 down to 103..221)
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new price-report record
$ws.Range("A102").Value = 7
$ws.Range("B102").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C102").Value = "Ñuble"
$ws.Range("D102").Value = 44664
$ws.Range("E102").Value = 16
$ws.Range("F102").Value = 100112003
$ws.Range("G102").Value = "Ajo"
$ws.Range("H102").Value = "Chino"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 100
$ws.Range("K102").Value = 20000
$ws.Range("L102").Value = 21000
$ws.Range("M102").Value = 20500
$ws.Range("N102").Value = "$/caja 10 kilos"
$ws.Range("O102").Value = "China"
$ws.Range("P102").Value = 2050
$ws.Range("Q102").Value = 10
$ws.Range("R102").Value = "Hortaliza"
